$d = $word.ActiveDocument

# 1) Three inline text replacements: "(Statement, X, Y, Z);" -> "(Statement / ParentKind, X, Y, Z);"
[void]$d.Content.Find.Execute("(Statement, SubjectKind, Predicate, Object);", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(Statement / ParentKind, SubjectKind, Predicate, Object);", 2)

[void]$d.Content.Find.Execute("(Statement, Subject, PredicateKind, Object);", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(Statement / ParentKind, Subject, PredicateKind, Object);", 2)

[void]$d.Content.Find.Execute("(Statement, Predicate, Subject, ObjectKind);", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(Statement / ParentKind, Predicate, Subject, ObjectKind);", 2)

$orderText = "Order: Kind interface. Kind / ParentKind hierarchical order relation: more abstract / more specific hierarchy tree nodes until singleton Kind."

# 2) Near "(Context, Statement, Kind, Resource);": the following empty paragraph gets the
#    "Order: Kind interface..." text, then two new empty paragraphs are inserted before "Augmentations:".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "(Context, Statement, Kind, Resource);") {
        $emptyPara = $d.Paragraphs.Item($i + 1)
        $emptyPara.Range.Text = $orderText

        $r = $d.Paragraphs.Item($i + 1).Range
        $r.InsertParagraphAfter()
        $r.InsertParagraphAfter()
        break
    }
}

# 3) After "Relationship: sameAs Statements.": insert a new empty paragraph, then a new
#    "Order: Kind interface..." paragraph (before the existing blank line that precedes
#    "Relationship Order / Comparison...").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "Relationship: sameAs Statements.") {
        $r = $p.Range
        $r.InsertParagraphAfter()
        $d.Paragraphs.Item($i + 1).Range.InsertParagraphAfter()
        $d.Paragraphs.Item($i + 2).Range.Text = $orderText
        break
    }
}

Write-Host "done"
